# Week2_Project_Design.pptx - "add W3 patterns to slide"
#
# Append a brand-new slide 6 ("Next Week") with links to the
# State and Observer design-pattern references.

$p = $ppt.ActivePresentation

# --- add the new "Next Week" slide (same Title+Content layout) --------
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)

# Title
$title = $newSlide.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Next Week"

# Body: StatePattern + link, then ObseverPattern + link
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "StatePattern"
[void]$body.InsertAfter(":")
[void]$body.InsertAfter("`rhttps://")
[void]$body.InsertAfter("en.wikipedia.org")
[void]$body.InsertAfter("/wiki/")
[void]$body.InsertAfter("State_pattern")
[void]$body.InsertAfter("`rObseverPattern")
[void]$body.InsertAfter(":")
[void]$body.InsertAfter("`rhttp://")
[void]$body.InsertAfter("www.oodesign.com")
[void]$body.InsertAfter("/observer-")
[void]$body.InsertAfter("pattern.html")

# The two URL paragraphs are plain (no bullet), like the source deck.
$body.Paragraphs(2, 1).ParagraphFormat.Bullet.Visible = 0
$body.Paragraphs(4, 1).ParagraphFormat.Bullet.Visible = 0
